# Allowed camp committee members to generate reports:
# - Corrected the "ADM Camp" slot count (EEE Camp row) from 90 to 89.
# - Added two new camp entries (rows 9 and 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix slot count for row 3 (EEE Camp) ---
$ws.Range("H3").Value = 89

# --- New row 9 ---
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "OOPs"
$ws.Range("C9").Value = 45286
$ws.Range("D9").Value = 45290
$ws.Range("E9").Value = 45285
$ws.Range("F9").Value = "SCSE"
$ws.Range("G9").Value = "Your MOMS' HOUSE"
$ws.Range("H9").Value = 8
$ws.Range("I9").Value = 10
$ws.Range("J9").Value = "rampant sex"
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = $true

# --- New row 10 ---
$ws.Range("A10").Value = 12
$ws.Range("B10").Value = "OOPsy"
$ws.Range("C10").Value = 45286
$ws.Range("D10").Value = 45290
$ws.Range("E10").Value = 45285
$ws.Range("F10").Value = "NTU"
$ws.Range("G10").Value = "your mums house"
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = "lololol"
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = $true

# --- Date formatting for the new rows' date columns (Start/End/Registration Deadline) ---
$ws.Range("C9:E10").Style = "Normal"
$ws.Range("C9:E10").NumberFormat = "dd/mm/yyyy"
